$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects_Toolbar")

# Row 3 ("Toolbar - Items"): switch the locator from an id to an xpath
$ws.Range("C3").Value = "xpath"
$ws.Range("E3").Value = ".//*[@id='Ribbon.ListItem-title']/a"

# Reflect the newly selected cell on the sheet (as recorded in the saved workbook)
$ws.Range("B3").Select()
